$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("harp expander v1")

# --- 1. Bump quantities / reference designators for parts shared with the new buffer IC ---

# 100nF 0402 caps: now also used as IC4 decoupling cap (C21)
$ws.Range("A5").Value2 = 11
$ws.Range("D5").Value2 = "C2, C6, C8, C9, C10, C12, C13, C14, C15, C16, C21"

# 1k 0402 resistors: now also used as IC4 pull-up/series resistor (R8)
$ws.Range("A9").Value2 = 5
$ws.Range("D9").Value2 = "R2, R3, R4, R6, R8"

# RED LED row: fix typo and add the new digital output tag LED (OUT0)
$ws.Range("A15").Value2 = 4
$ws.Range("D15").Value2 = "TAG_B, TAG_T, STATE, OUT0"

# --- 2. R7 (39 ohm resistor) footprint is now confirmed: fill package, drop the "Not sure yet" note ---
$ws.Range("C8").Value2 = "R0402"
$ws.Range("I8").ClearContents()
$ws.Range("E8:H8").Font.Color = $ws.Range("D8").Font.Color

# --- 3. Insert the new IC4 digital-output buffer row (row 20, pushes OSC1.. down by one) ---
$ws.Rows("20:20").Insert()

$ws.Range("A20").Value2 = 1
$ws.Range("B20").Value2 = "IC BUFFER NON-INVERT 5.5V"
$ws.Range("C20").Value2 = "TSOP-5"
$ws.Range("D20").Value2 = "IC4"
$ws.Range("E20").Value2 = "M74VHC1GT126DT1G"
$ws.Range("F20").Value2 = "SMD"
$ws.Range("G20").Value2 = "Digikey"
$ws.Range("H20").Value2 = "M74VHC1GT126DT1GOSDKR-ND"

# --- 4. Refresh the summary counters (unique part count, SMD part count) ---
$ws.Range("F29").Value2 = 21
$ws.Range("F30").Value2 = 41

# --- 5. Cost table: the "components" line now also covers the extra cap/resistor/buffer ---
$ws.Range("F35").Formula = "=49.74+0.34+1.226"

# --- 6. Match the saved selection from the edit ---
$ws.Range("E21:H21").Select()
